$d = $word.ActiveDocument

# Locate the anchor paragraph (the existing last bullet, 'In `Match` table...')
# and append a placeholder paragraph after it, inheriting the same
# ListParagraph / numPr (bulleted list) paragraph formatting.
$old = "In ‘Match’ table, the very last columns are home win odds, draw odds and away win odds computed using different methods. But all these seem to be giving results in the same ratio. We could take the average of these numbers and reduce it to just three columns, “Home win odds”, “Draw Odds”, and “Away win odds”."
$new = "In ‘Match’ table, the very last columns are home win odds, draw odds and away win odds computed using different methods. But all these seem to be giving results in the same ratio. We could take the average of these numbers and reduce it to just three columns, “Home win odds”, “Draw Odds”, and “Away win odds”.^p@@NEWPARA1@@"
$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) {
    throw "Could not locate the anchor paragraph to extend."
}

function Get-LastParagraph {
    return $d.Paragraphs($d.Paragraphs.Count)
}

# Fill in the first new paragraph. It was originally typed as one
# sentence, then ' (183978 rows in player_Attributes table)' was
# inserted mid-sentence afterwards, so it ends up as three runs.
$p = Get-LastParagraph
$rStart = $p.Range.Start
$rEnd = $p.Range.End - 1
$r = $d.Range($rStart, $rEnd)
$r.Text = ""
$r.Collapse(1)
$r.InsertAfter("How are the player attributes stored in the tables? Also, what are the relationships of these attributes to the match? Can we use it to train our tree?")
$r.Collapse(0)
$r.InsertAfter(" (183978 rows in player_Attributes table)")
$r.Collapse(0)
$r.InsertAfter(". It has a lot of attributes, but the question is how to associate the attributes of a particular player to the match? Shall we consider the attributes of the team as the aggregation of the attributes of the individual player playing that match?")

# Append new paragraph 2 of 6 (same list formatting carried over).
$lastP = Get-LastParagraph
$endRng = $lastP.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$newP = Get-LastParagraph
$newRng = $newP.Range
$newRng.Collapse(0)
$newRng.InsertAfter("Also, there is only so much information given. For instance, the values of particular attributes are given but it is not provided how to interpret that information! ")

# Append new paragraph 3 of 6 (same list formatting carried over).
$lastP = Get-LastParagraph
$endRng = $lastP.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$newP = Get-LastParagraph
$newRng = $newP.Range
$newRng.Collapse(0)
$newRng.InsertAfter("Given an attribute, potential, it not provided how the attribute is measured. What is the scale of the attribute? Is it evaluated on a scale of 100 or it is indefinite? Is a higher number better or worse? How to compare two numbers of this attribute? Is it linear? Can we just subtract the number to compare which player has how much more potential compare to a particular player?")

# Append new paragraph 4 of 6 (same list formatting carried over).
$lastP = Get-LastParagraph
$endRng = $lastP.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$newP = Get-LastParagraph
$newRng = $newP.Range
$newRng.Collapse(0)
$newRng.InsertAfter("Also, while making a complex Decision Tree, should we give different weight to different attributes? How do we decide how much weight to give to which attribute? One way is to just make different trees with different weights assigned to different attributes, and pick the one that gives the best results.")

# Append new paragraph 5 of 6 (same list formatting carried over).
$lastP = Get-LastParagraph
$endRng = $lastP.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$newP = Get-LastParagraph
$newRng = $newP.Range
$newRng.Collapse(0)
$newRng.InsertAfter("One problem with this approach is we may never find the best decision tree also there are exponentially many trees possible with different weight attributes, we can’t compute them all.")

# Append new paragraph 6 of 6 (same list formatting carried over).
$lastP = Get-LastParagraph
$endRng = $lastP.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$newP = Get-LastParagraph
$newRng = $newP.Range
$newRng.Collapse(0)
$newRng.InsertAfter("Hit and miss without anything also gives 50% ! :P")

